# Update the single data row (row 2) in Hoja1 with the new delegacion record.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A2").Value = "MARCELINO UGARTE"
$ws.Range("B2").Value = "5"
$ws.Range("C2").Value = "99906"
$ws.Range("D2").Value = "9 DE JULIO"
$ws.Range("E2").Value = "9 DE JULIO"

# Move the active selection to B6, matching the saved view state.
$ws.Range("B6").Select()
